$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "277.56"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "0.96%"
$c = $ws.Range("G2")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "27.26"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "2.04%"
$c = $ws.Range("G3")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "4.873"
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "-0.28%"
$c = $ws.Range("G4")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "0.06417"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "1.28%"
$c = $ws.Range("G5")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "6.957"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "1.11%"
$c = $ws.Range("G6")
$c.NumberFormat = "@"
$c.Value = "8"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.8858"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "1.88%"
$c = $ws.Range("G7")
$c.NumberFormat = "@"
$c.Value = "8"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "1.180"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "-6.37%"
$c = $ws.Range("G8")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.1539"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "-2.28%"
$c = $ws.Range("G9")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.05144"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "-1.08%"
$c = $ws.Range("G10")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07423"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "-0.06%"
$c = $ws.Range("G11")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.02888"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "-1.59%"
$c = $ws.Range("G12")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.08965"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "-0.81%"
$c = $ws.Range("G13")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.001577"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "-0.16%"
$c = $ws.Range("G14")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.0006354"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "0.35%"
$c = $ws.Range("G15")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.006098"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "1.31%"
$c = $ws.Range("G16")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "3.481"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "0.89%"
$c = $ws.Range("G17")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "3.311"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "-0.22%"
$c = $ws.Range("G18")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "0.07%"
$c = $ws.Range("G19")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.3149"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "1.15%"
$c = $ws.Range("G20")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "0.97%"
$c = $ws.Range("G21")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "3.902"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "-0.09%"
$c = $ws.Range("G22")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.04421"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "1.23%"
$c = $ws.Range("G23")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "8.68%"
$c = $ws.Range("G24")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("G25")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.001177"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "-0.50%"
$c = $ws.Range("G26")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.003867"
$c = $ws.Range("G27")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "8.12%"
$c = $ws.Range("G28")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "15.41%"
$c = $ws.Range("G29")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("G30")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("G31")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("G32")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("G33")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("G34")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("G35")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("G36")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("G37")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("G38")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("G39")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.04162"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "1.45%"
$c = $ws.Range("G40")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.006771"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "-0.06%"
$c = $ws.Range("G41")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "0.60%"
$c = $ws.Range("G42")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.002010"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "-6.25%"
$c = $ws.Range("G43")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.01148"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "6.62%"
$c = $ws.Range("G44")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.00005302"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "-0.21%"
$c = $ws.Range("G45")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("G46")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "-11.95%"
$c = $ws.Range("G47")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("G48")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("G49")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("G50")
$c.NumberFormat = "@"
$c.Value = "8"
$c = $ws.Range("G51")
$c.NumberFormat = "@"
$c.Value = "8"
